# Fruta / hortaliza, semanal
# Insert two new weekly price rows (2023-12-11, "Provincia de Limarí") at the top
# of the data block (rows 50-51), shifting the existing rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 50; existing rows 50:130 shift to 52:132.
$ws.Rows("50:51").Insert()

# --- New row 50 ---
$ws.Range("A50").Value = 8
$ws.Range("B50").Value = "Terminal La Palmera de La Serena"
$ws.Range("C50").Value = "Coquimbo"
$ws.Range("D50").Value = 45271
$ws.Range("E50").Value = 4
$ws.Range("F50").Value = 100112027
$ws.Range("G50").Value = "Melón"
$ws.Range("H50").Value = "Tuna"
$ws.Range("I50").Value = "Primera"
$ws.Range("J50").Value = 2400
$ws.Range("K50").Value = 1900
$ws.Range("L50").Value = 2000
$ws.Range("M50").Value = 1950
$ws.Range("N50").Value = "`$/unidad"
$ws.Range("O50").Value = "Provincia de Limarí"
$ws.Range("P50").Value = 1950
$ws.Range("Q50").Value = 1
$ws.Range("R50").Value = "Hortaliza"

# --- New row 51 ---
$ws.Range("A51").Value = 8
$ws.Range("B51").Value = "Terminal La Palmera de La Serena"
$ws.Range("C51").Value = "Coquimbo"
$ws.Range("D51").Value = 45271
$ws.Range("E51").Value = 4
$ws.Range("F51").Value = 100112027
$ws.Range("G51").Value = "Melón"
$ws.Range("H51").Value = "Tuna"
$ws.Range("I51").Value = "Segunda"
$ws.Range("J51").Value = 3600
$ws.Range("K51").Value = 1400
$ws.Range("L51").Value = 1500
$ws.Range("M51").Value = 1450
$ws.Range("N51").Value = "`$/unidad"
$ws.Range("O51").Value = "Provincia de Limarí"
$ws.Range("P51").Value = 1450
$ws.Range("Q51").Value = 1
$ws.Range("R51").Value = "Hortaliza"

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()
